$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.933.36"
$ws.Range("E2").Value = "  -5.72%  "

# Row 3
$ws.Range("D3").Value = "3.738.70"
$ws.Range("E3").Value = "  -6.26%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.97%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.645"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.93%  "

# Row 8
$ws.Range("E8").Value = "  +0.51%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.15%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.74%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.27%  "

# Row 12
$ws.Range("E12").Value = "  -6.63%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.48%  "

# Row 14
$ws.Range("D14").Value = "4.344.23"
$ws.Range("E14").Value = "  -5.99%  "

# Row 15
$ws.Range("D15").Value = "3.769.80"
$ws.Range("E15").Value = "  -5.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.07%  "

# Row 18
$ws.Range("E18").Value = "  -8.13%  "

# Row 19
$ws.Range("E19").Value = "  -2.73%  "

# Row 20
$ws.Range("D20").Value = "68.833.57"
$ws.Range("E20").Value = "  -5.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "416.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "90.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.89%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.20%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.89%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.47%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.64%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "45.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.71%  "

# Row 34
$ws.Range("E34").Value = "  -7.45%  "

# Row 35
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.38%  "

# Row 36
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0949"
$ws.Range("E36").Value = "  -6.74%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "608.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.90%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.410"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.84%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.141"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.13%  "

# Row 44
$ws.Range("E44").Value = "  -7.66%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.81%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -12.53%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.138"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.57%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.94%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.786.83"
$ws.Range("E49").Value = "  -3.12%  "

# Row 50
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000273"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.59%  "

# Row 51
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -19.19%  "
